$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2

        if ($val -is [string]) {
            $newVal = $val
            $newVal = $newVal.Replace("D64", "D69")
            $newVal = $newVal.Replace("D80", "D86")
            $newVal = $newVal.Replace("D51", "D55")
            $newVal = $newVal.Replace("S30", "S31")

            if ($newVal -ne $val) {
                $cell.Value2 = $newVal
            }
        }
    }
}
